$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  "words",
  "adalat",
  "aspirin",
  "atorva",
  "atozet",
  "banan",
  "calcio",
  "cozaar",
  "dilatrend",
  "doxycycline",
  "Edoxaban",
  "hemonia",
  "lasix",
  "laxative",
  "MTM",
  "neotigasone",
  "nephvita",
  "NSAID",
  "NSAIDs",
  "pantoloc",
  "plavix",
  "renvela",
  "torem",
  "torsemide",
  "valtrex",
  "그리세틸",
  "글리아타민",
  "글리킨연질",
  "글리피딘정",
  "글립타이드정",
  "네비레트정",
  "네시나액트정",
  "넥리카캡슐",
  "넥사론정",
  "넥시움",
  "넥포정",
  "니페딕스",
  "다비듀오",
  "데잘탄정",
  "덱실란트",
  "덱실란트디알캡슐",
  "디아미크롱",
  "디오르탄정",
  "디오반",
  "디카맥스",
  "라바졸정",
  "라베듀오",
  "라베스란정",
  "라식스",
  "라식스",
  "라식스정",
  "레나메진",
  "레날민정",
  "레날빈",
  "레미코트",
  "레미피드정",
  "렌벨라",
  "로바젯정",
  "로수젯정",
  "로포타현탁액",
  "리나치올",
  "리바로젯",
  "리피딜슈프라정",
  "리피토",
  "릭시아나",
  "마프론정",
  "바로디핀",
  "바스티난",
  "벤조다이아제핀",
  "브루펜",
  "삐콤정",
  "삭센다",
  "세레브렉스",
  "세보칸정",
  "세비카정",
  "시나롱",
  "실로스탄씨알정",
  "씨엠피나정",
  "아낙성",
  "아다핀오스모서방정",
  "아달라트",
  "아마릴",
  "아모딘정",
  "아모스타정",
  "아서틸",
  "아세라노정",
  "아젯듀오",
  "아토넷정",
  "아토르바스타틴",
  "아토르바스타틴정",
  "아토르바정",
  "아토칸정",
  "아트로정",
  "안탁스",
  "알레세틴정",
  "암로디핀정",
  "엔테론",
  "엔테론정",
  "엔트레스토",
  "엘리퀴스",
  "오마코",
  "올메텍",
  "울트라셋",
  "울프라셑",
  "유로박솜",
  "이가탄",
  "인벨라정",
  "자누매트정",
  "자누비아",
  "자니딥정",
  "자디앙",
  "자렐토",
  "제미글로정",
  "직듀오",
  "카리메이트",
  "카리메트산",
  "카비드",
  "칸다벨정",
  "케이슈머산",
  "케이캡",
  "코다론",
  "코푸",
  "큐미날",
  "크레더블",
  "크레메진",
  "클로스원",
  "타스나",
  "테넬리아",
  "텔미플러스",
  "토르셈정",
  "트라몰",
  "트라젠타정",
  "팔로델",
  "페브락",
  "페브릭",
  "페브릭정",
  "펜잘",
  "펜폴",
  "포스레놀산",
  "포스포마이신",
  "프레탈",
  "프로맥정",
  "플라빅스",
  "하이손정",
  "하이페질정",
  "항고혈압제",
  "헤르벤",
  "훼로바",
  "훼로바유"
)

for ($i = 0; $i -lt $values.Count; $i++) {
  $ws.Cells.Item($i + 1, 1).Value2 = $values[$i]
}
